# This sheet is a weekly price log for "Locoto" (Vega Modelo de Temuco).
# The commit adds one new week's record at the top of the data block
# (row 49), pushing the existing rows 49-83 down by one (to 50-84).
#
# Insert a new row at row 49 (shifts rows 49:83 down to 50:84, including
# all formatting/styles carried on those rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("49:49").Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Cells.Item(49, 1).Value = 10
$ws.Cells.Item(49, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(49, 3).Value = "La Araucanía"
$ws.Cells.Item(49, 4).Value = 45086
$ws.Cells.Item(49, 5).Value = 9
$ws.Cells.Item(49, 6).Value = 100112042
$ws.Cells.Item(49, 7).Value = "Locoto"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 80
$ws.Cells.Item(49, 11).Value = 4400
$ws.Cells.Item(49, 12).Value = 4400
$ws.Cells.Item(49, 13).Value = 4400
$ws.Cells.Item(49, 14).Value = "$/kilo"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 4400
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# Make sure the D column on the new row keeps the same date/time display
# style ("s=2" / numFmtId 165) that every other row in this block uses.
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
